$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Locate the "A step by step guide ..." paragraph (EuroScope Setup /
# Observing the network guide links) without hard-coding its index.
# ---------------------------------------------------------------------
$targetIndex = -1
$paras = $d.Paragraphs
$count = $paras.Count
for ($i = 1; $i -le $count; $i++) {
    $candidate = $paras.Item($i)
    if ($candidate.Range.Text.StartsWith("A step by step guide on setting EuroScope")) {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate the EuroScope step-by-step guide paragraph"
}

# ---------------------------------------------------------------------
# Replace the whole paragraph body with the restructured runs /
# hyperlinks. The two hyperlinks keep their original relationship ids
# (rId11 -> EuroScope setup guide, rId12 -> observing/forum guide) so
# their target URLs are preserved untouched.
# ---------------------------------------------------------------------
$paraRange = $paras.Item($targetIndex).Range

$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">' +
    '<w:r><w:t>S</w:t></w:r>' +
    '<w:r><w:t>tep by step</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> guides for </w:t></w:r>' +
    '<w:hyperlink r:id="rId11" w:history="1">' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>EuroScope</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> setup</w:t></w:r>' +
    '</w:hyperlink>' +
    '<w:r><w:t xml:space="preserve"> &amp; </w:t></w:r>' +
    '<w:hyperlink r:id="rId12" w:history="1">' +
        '<w:r><w:t>Observing the network</w:t></w:r>' +
    '</w:hyperlink>' +
    '<w:r><w:t xml:space="preserve">, alongside frequently asked questions are available on the </w:t></w:r>' +
    '<w:r><w:t>VATSIM UK Documentation Site</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '</w:p>'

$paraRange.InsertXML($newParaXml) | Out-Null

# ---------------------------------------------------------------------
# InsertXML drops <w:rStyle> references (character-style links), so
# re-apply the "Hyperlink" character style to the display text of each
# hyperlink run now that the plain text/structure is in place.
# ---------------------------------------------------------------------
$paraRange = $paras.Item($targetIndex).Range

$rng = $d.Range($paraRange.Start, $paraRange.End)
$rng.Find.Execute("EuroScope", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Style = "Hyperlink"

$paraRange = $paras.Item($targetIndex).Range
$rng = $d.Range($paraRange.Start, $paraRange.End)
$rng.Find.Execute(" setup", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Style = "Hyperlink"

$paraRange = $paras.Item($targetIndex).Range
$rng = $d.Range($paraRange.Start, $paraRange.End)
$rng.Find.Execute("Observing the network", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Style = "Hyperlink"
